# CIERRE 9 oCT 2021
#
# Payroll ("RECIBOS NOMINA 2020") closing update on Hoja1:
#   - Quincena (pay period) 1: días trabajados 5->6, percepción 2166->2600,
#     deducción 0->433 (totals recalc via existing SUM formulas).
#   - Quincena 2: percepción 0->3080 (total recalcs via SUM formula).
#   - Quincena 3: deducción 0->1250 (total recalcs via SUM formula).
#   - The four TODAY()/chained-date cells recalc automatically to the new
#     closing date once the workbook is recalculated.
#   - Leaves the view scrolled/selected on the last period's total (K41),
#     matching where the user finished editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # Hoja1 is the (only) selected/active sheet

# --- Quincena 1 (rows 3-6) ---
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 2600
$ws.Range("K4").Value = 433

# --- Quincena 2 (rows 21-24) ---
$ws.Range("K21").Value = 3080

# --- Quincena 3 (rows 38-41) ---
$ws.Range("K40").Value = 1250

# Recalculate so the SUM() totals and the TODAY()-chained date cells
# (C14/I14/C32/I32/C48/I48/C65) pick up fresh cached values.
$excel.Calculate()

# --- Leave the sheet scrolled to / selection on the last edited total ---
$ws.Activate()
$excel.Goto($ws.Range("K41"), $true)
$ws.Range("K41").Select()
